# This workbook lists outbreak clusters with their active-case counts and
# whether the row is "new" or "old" (rows 1:C header + 43 data rows -> A1:C44).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear existing contents (keeps formatting/styles) before laying down the
# fully refreshed table in a single write.
$ws.Cells.ClearContents()

$data = New-Object 'object[,]' 44,3
$data[0,0] = "Cluster name"
$data[0,1] = "Active cases"
$data[0,2] = "Exist"
$data[1,0] = '139 Highett St Apartment Complex Richmond'
$data[1,1] = 10
$data[1,2] = 'new'
$data[2,0] = '139 Highett St Apartment Complex Richmond'
$data[2,1] = 11
$data[2,2] = 'old'
$data[3,0] = '3175 The Bays Aged Care Facility Hastings'
$data[3,1] = 16
$data[3,2] = 'old'
$data[4,0] = '3175 The Bays Aged Care Facility Hastings'
$data[4,1] = 17
$data[4,2] = 'new'
$data[5,0] = '3612 BlueCross Glengowrie Outbreak'
$data[5,1] = 30
$data[5,2] = 'old'
$data[6,0] = '3612 BlueCross Glengowrie Outbreak'
$data[6,1] = 33
$data[6,2] = 'new'
$data[7,0] = '3684 Homestyle Aged Care Langford Grange Cranbourne East Outbreak'
$data[7,1] = 20
$data[7,2] = 'old'
$data[8,0] = '3684 Homestyle Aged Care Langford Grange Cranbourne East Outbreak'
$data[8,1] = 23
$data[8,2] = 'new'
$data[9,0] = '4075 Ferndale Gardens Aged Care Services Bayswater North Outbreak'
$data[9,1] = 16
$data[9,2] = 'old'
$data[10,0] = '4075 Ferndale Gardens Aged Care Services Bayswater North Outbreak'
$data[10,1] = 17
$data[10,2] = 'new'
$data[11,0] = '4518 Regis Aged Care Fawkner Outbreak'
$data[11,1] = 13
$data[11,2] = 'new'
$data[12,0] = 'AW Window Transport Group Depot North Geelong Outbreak'
$data[12,1] = 10
$data[12,2] = 'new'
$data[13,0] = 'Australian Lamb Colac East'
$data[13,1] = 13
$data[13,2] = 'old'
$data[14,0] = 'Bread Solutions Braeside Outbreak'
$data[14,1] = 18
$data[14,2] = 'new'
$data[15,0] = 'Bread Solutions Braeside Outbreak'
$data[15,1] = 19
$data[15,2] = 'old'
$data[16,0] = 'CS Square Caroline Springs Outbreak'
$data[16,1] = 14
$data[16,2] = 'new'
$data[17,0] = 'CS Square Caroline Springs Outbreak'
$data[17,1] = 17
$data[17,2] = 'old'
$data[18,0] = 'Cedar Meats Australia Brooklyn Outbreak'
$data[18,1] = 10
$data[18,2] = 'new'
$data[19,0] = 'Cedar Meats Australia Brooklyn Outbreak'
$data[19,1] = 11
$data[19,2] = 'old'
$data[20,0] = 'Community Kids Pascoe Vale Early Education Centre Pascoe Vale Outbreak'
$data[20,1] = 20
$data[20,2] = 'new'
$data[21,0] = 'Embracia Aged Care Reservoir Outbreak'
$data[21,1] = 19
$data[21,2] = 'new'
$data[22,0] = 'Embracia Aged Care Reservoir Outbreak'
$data[22,1] = 22
$data[22,2] = 'old'
$data[23,0] = 'Guardian Childcare Caulfield Outbreak'
$data[23,1] = 17
$data[23,2] = 'old'
$data[24,0] = 'Guardian Childcare Caulfield Outbreak'
$data[24,1] = 20
$data[24,2] = 'new'
$data[25,0] = 'Inghams Enterprise Somerville Outbreak'
$data[25,1] = 15
$data[25,2] = 'old'
$data[26,0] = 'Inghams Enterprise Somerville Outbreak'
$data[26,1] = 17
$data[26,2] = 'new'
$data[27,0] = 'Launch Housing City Edge Crisis Accommodation South Melbourne'
$data[27,1] = 11
$data[27,2] = 'new'
$data[28,0] = 'Launch Housing City Edge Crisis Accommodation South Melbourne'
$data[28,1] = 12
$data[28,2] = 'old'
$data[29,0] = 'Northern Health Northern Hospital Epping Emergency Department Tier 1B'
$data[29,1] = 41
$data[29,2] = 'new'
$data[30,0] = 'Northern Health Northern Hospital Epping Emergency Department Tier 1B'
$data[30,1] = 42
$data[30,2] = 'old'
$data[31,0] = 'Northern Health The Northern Hospital Epping'
$data[31,1] = 14
$data[31,2] = 'new'
$data[32,0] = 'Northern Health The Northern Hospital Epping'
$data[32,1] = 15
$data[32,2] = 'old'
$data[33,0] = 'Robin Hood Inn Drouin West Outbreak'
$data[33,1] = 41
$data[33,2] = 'old'
$data[34,0] = 'Robin Hood Inn Drouin West Outbreak'
$data[34,1] = 42
$data[34,2] = 'new'
$data[35,0] = 'Social Gathering Warrnambool 28 Sep Outbreak'
$data[35,1] = 13
$data[35,2] = 'old'
$data[36,0] = 'St Vincents Hospital Emergency Department Melbourne'
$data[36,1] = 41
$data[36,2] = 'new'
$data[37,0] = 'St Vincents Hospital Emergency Department Melbourne'
$data[37,1] = 42
$data[37,2] = 'old'
$data[38,0] = 'Target Distribution Centre Truganina Outbreak'
$data[38,1] = 19
$data[38,2] = 'new'
$data[39,0] = 'Target Distribution Centre Truganina Outbreak'
$data[39,1] = 20
$data[39,2] = 'old'
$data[40,0] = 'Visy Recycling Springvale'
$data[40,1] = 21
$data[40,2] = 'new'
$data[41,0] = 'Visy Recycling Springvale'
$data[41,1] = 29
$data[41,2] = 'old'
$data[42,0] = 'Werribee Mercy Hospital Emergency Department'
$data[42,1] = 24
$data[42,2] = 'old'
$data[43,0] = 'Werribee Mercy Hospital Emergency Department'
$data[43,1] = 25
$data[43,2] = 'new'

$ws.Range("A1:C44").Value = $data
